# Update "想去人数" (people-interested count) figures in column F across
# the 展览 (Exhibitions), 演出 (Shows) and 全部类型 (All types) sheets,
# matching the regenerated site data (gh-pages output @ 456a3b4).
# 本地生活 (Local life) only has a header row, so it needs no updates.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 20875
$ws.Range("F7").Value = 7803
$ws.Range("F9").Value = 22
$ws.Range("F11").Value = 297
$ws.Range("F13").Value = 181
$ws.Range("F15").Value = 27
$ws.Range("F17").Value = 215
$ws.Range("F19").Value = 502
$ws.Range("F21").Value = 701
$ws.Range("F27").Value = 46
$ws.Range("F29").Value = 208
$ws.Range("F31").Value = 592
$ws.Range("F32").Value = 115
$ws.Range("F33").Value = 4972
$ws.Range("F36").Value = 28
$ws.Range("F38").Value = 12943
$ws.Range("F40").Value = 112
$ws.Range("F41").Value = 44
$ws.Range("F42").Value = 66
$ws.Range("F43").Value = 295
$ws.Range("F44").Value = 408
$ws.Range("F45").Value = 4040
$ws.Range("F46").Value = 325

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 320

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 20875
$ws.Range("F7").Value = 7803
$ws.Range("F9").Value = 22
$ws.Range("F11").Value = 297
$ws.Range("F13").Value = 181
$ws.Range("F15").Value = 27
$ws.Range("F17").Value = 215
$ws.Range("F19").Value = 502
$ws.Range("F21").Value = 701
$ws.Range("F27").Value = 46
$ws.Range("F29").Value = 208
$ws.Range("F30").Value = 320
$ws.Range("F31").Value = 592
$ws.Range("F33").Value = 115
$ws.Range("F35").Value = 4972
$ws.Range("F38").Value = 28
$ws.Range("F40").Value = 12943
$ws.Range("F42").Value = 112
$ws.Range("F43").Value = 44
$ws.Range("F44").Value = 66
$ws.Range("F45").Value = 295
$ws.Range("F46").Value = 408
$ws.Range("F47").Value = 4040
$ws.Range("F48").Value = 325
